$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for columns A..T (1..20)
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# New full data for rows 2..13 (20 columns: A..T).
# Rows 2-7: sending cluster "ECs" paired with receiving clusters ECs/FAPs/Inflammatory-Mac/MuSCs/Neutrophils/Resolving-Mac
# Rows 8-13: sending cluster "FAPs" (same as original rows 2-7, target clusters) re-normalised now that ECs is also a sender
$rowsData = @(
    @("ECs", "Rspo3", "Sdc4", "ECs", 1, 0.3333333333333333, 0.01818866666666667, 0.054566, 0.006403810693375696, 0.006403810693375696, 3, 1, 3.264166666666667, 9.7925, 0.02701285365228691, 0.02701285365228691, 0.05937083944444446, 0.534337555, 0.0001729852010771076, 0.0001729852010771076),
    @("ECs", "Rspo3", "Sdc4", "FAPs", 1, 0.3333333333333333, 0.01818866666666667, 0.054566, 0.006403810693375696, 0.006403810693375696, 3, 1, 37.12743, 111.38229, 0.3072508041078968, 0.3072508041078968, 0.67529844846, 6.07768603614, 0.001967575984894431, 0.001967575984894431),
    @("ECs", "Rspo3", "Sdc4", "Inflammatory-Mac", 1, 0.3333333333333333, 0.01818866666666667, 0.054566, 0.006403810693375696, 0.006403810693375696, 3, 1, 15.075014, 45.225042, 0.1247543978518794, 0.1247543978518794, 0.2741944046413334, 2.467749641772, 0.0007989035470095113, 0.0007989035470095113),
    @("ECs", "Rspo3", "Sdc4", "MuSCs", 1, 0.3333333333333333, 0.01818866666666667, 0.054566, 0.006403810693375696, 0.006403810693375696, 3, 1, 31.08386466666667, 93.25159400000001, 0.2572368303869774, 0.2572368303869774, 0.5653740531337779, 5.088366478204001, 0.001647295965162196, 0.001647295965162196),
    @("ECs", "Rspo3", "Sdc4", "Neutrophils", 1, 0.3333333333333333, 0.01818866666666667, 0.054566, 0.006403810693375696, 0.006403810693375696, 3, 1, 5.556591666666667, 16.669775, 0.04598398697896871, 0.04598398697896871, 0.1010669936277778, 0.9096029426500001, 0.0002944727475399686, 0.0002944727475399686),
    @("ECs", "Rspo3", "Sdc4", "Resolving-Mac", 1, 0.3333333333333333, 0.01818866666666667, 0.054566, 0.006403810693375696, 0.006403810693375696, 3, 1, 28.73046866666667, 86.191406, 0.2377611270219907, 0.2377611270219907, 0.5225689177551112, 4.703120259796, 0.001522577247692481, 0.001522577247692481),
    @("FAPs", "Rspo3", "Sdc4", "ECs", 3, 1, 2.822099333333334, 8.466298, 0.9935961893066243, 0.9935961893066244, 3, 1, 3.264166666666667, 9.7925, 0.02701285365228691, 0.02701285365228691, 9.211802573888889, 82.906223165, 0.0268398684512098, 0.02683986845120981),
    @("FAPs", "Rspo3", "Sdc4", "FAPs", 3, 1, 2.822099333333334, 8.466298, 0.9935961893066243, 0.9935961893066244, 3, 1, 37.12743, 111.38229, 0.3072508041078968, 0.3072508041078968, 104.77729545138, 942.99565906242, 0.3052832281230024, 0.3052832281230025),
    @("FAPs", "Rspo3", "Sdc4", "Inflammatory-Mac", 3, 1, 2.822099333333334, 8.466298, 0.9935961893066243, 0.9935961893066244, 3, 1, 15.075014, 45.225042, 0.1247543978518794, 0.1247543978518794, 42.54318695939067, 382.888682634516, 0.1239554943048699, 0.1239554943048699),
    @("FAPs", "Rspo3", "Sdc4", "MuSCs", 3, 1, 2.822099333333334, 8.466298, 0.9935961893066243, 0.9935961893066244, 3, 1, 31.08386466666667, 93.25159400000001, 0.2572368303869774, 0.2572368303869774, 87.72175375322357, 789.4957837790121, 0.2555895344218152, 0.2555895344218153),
    @("FAPs", "Rspo3", "Sdc4", "Neutrophils", 3, 1, 2.822099333333334, 8.466298, 0.9935961893066243, 0.9935961893066244, 3, 1, 5.556591666666667, 16.669775, 0.04598398697896871, 0.04598398697896871, 15.68125363810556, 141.13128274295, 0.04568951423142874, 0.04568951423142874),
    @("FAPs", "Rspo3", "Sdc4", "Resolving-Mac", 3, 1, 2.822099333333334, 8.466298, 0.9935961893066243, 0.9935961893066244, 3, 1, 28.73046866666667, 86.191406, 0.2377611270219907, 0.2377611270219907, 81.08023647055423, 729.722128234988, 0.2362385497742982, 0.2362385497742983)
)

for ($i = 0; $i -lt $rowsData.Length; $i++) {
    $r = $i + 2
    $rowVals = $rowsData[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}
